$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename columns
#   A1: "No. Nómina" -> "NumNomina"
#   B1: "Nombre(s)"  -> "Nombres"
# (C1 "Apellidos" stays the same)
$ws.Range("A1").Value = "NumNomina"
$ws.Range("B1").Value = "Nombres"

# Move the active selection to D14 (cosmetic cursor-position change
# left behind in the saved file)
$ws.Range("D14").Select()
